# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# to reflect the newly generated report times.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: G2 "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-05 19:12:40"

# --- zh-cn sheet: H2 "Correspond Handoff Datetime", K2 "Correspond Handback DateTime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-05 19:12:35"
$wsZhCn.Range("K2").Value = "2016-09-05 19:13:01"

# --- de-de sheet: H2 "Correspond Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-05 19:13:15"
